# Apply translation-related header renames and view/selection changes.

$wb = $excel.ActiveWorkbook

# --- settings sheet: rename display.title -> display.title.text ---
# (done first so the new shared string is appended before the survey ones,
#  matching the author's original edit order)
$settings = $wb.Worksheets.Item("settings")
$settings.Range("C1").Value = "display.title.text"
$settings.Range("C2").Select()

# --- survey sheet: rename display column headers ---
$survey = $wb.Worksheets.Item("survey")
$survey.Range("H1").Value = "display.prompt.text"
$survey.Range("I1").Value = "display.hint.text"
$survey.Range("K1").Value = "display.button_label.text"

# widen columns J and K on the survey sheet
$survey.Columns.Item(10).ColumnWidth = 9.33203125
$survey.Columns.Item(11).ColumnWidth = 22.33203125

# scroll/select state on survey sheet
$survey.Range("C1").Select()
$survey.Application.ActiveWindow.ScrollColumn = 3
$survey.Range("K23").Select()

# --- properties sheet selection unchanged, just no longer the active tab ---
$properties = $wb.Worksheets.Item("properties")
$properties.Range("F13").Select()

# --- choices sheet becomes the active tab ---
$choices = $wb.Worksheets.Item("choices")
$choices.Activate()
$choices.Range("C20").Select()
